$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = '="70.719.58"'
$ws.Range('D2').Copy()
$ws.Range('D2').PasteSpecial(-4163)
$ws.Range('E2').Formula = '="  +2.81%  "'
$ws.Range('E2').Copy()
$ws.Range('E2').PasteSpecial(-4163)
$ws.Range('D3').Formula = '="3.572.49"'
$ws.Range('D3').Copy()
$ws.Range('D3').PasteSpecial(-4163)
$ws.Range('E3').Formula = '="  +2.02%  "'
$ws.Range('E3').Copy()
$ws.Range('E3').PasteSpecial(-4163)
$ws.Range('E4').Formula = '="  -0.04%  "'
$ws.Range('E4').Copy()
$ws.Range('E4').PasteSpecial(-4163)
$ws.Range('D5').Formula = '="599.97"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Formula = '="  +2.38%  "'
$ws.Range('E5').Copy()
$ws.Range('E5').PasteSpecial(-4163)
$ws.Range('D6').Formula = '="173.25"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Formula = '="  +2.20%  "'
$ws.Range('E6').Copy()
$ws.Range('E6').PasteSpecial(-4163)
$ws.Range('D7').Formula = '="3.566.47"'
$ws.Range('D7').Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E8').Formula = '="  +1.26%  "'
$ws.Range('E8').Copy()
$ws.Range('E8').PasteSpecial(-4163)
$ws.Range('E9').Formula = '="  +0.01%  "'
$ws.Range('E9').Copy()
$ws.Range('E9').PasteSpecial(-4163)
$ws.Range('E10').Formula = '="  +6.73%  "'
$ws.Range('E10').Copy()
$ws.Range('E10').PasteSpecial(-4163)
$ws.Range('D11').Formula = '="7.39"'
$ws.Range('D11').Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Formula = '="  +9.35%  "'
$ws.Range('E11').Copy()
$ws.Range('E11').PasteSpecial(-4163)
$ws.Range('E12').Formula = '="  +2.67%  "'
$ws.Range('E12').Copy()
$ws.Range('E12').PasteSpecial(-4163)
$ws.Range('D13').Formula = '="46.69"'
$ws.Range('D13').Copy()
$ws.Range('D13').PasteSpecial(-4163)
$ws.Range('E13').Formula = '="  -0.21%  "'
$ws.Range('E13').Copy()
$ws.Range('E13').PasteSpecial(-4163)
$ws.Range('E14').Formula = '="  +1.70%  "'
$ws.Range('E14').Copy()
$ws.Range('E14').PasteSpecial(-4163)
$ws.Range('D15').Formula = '="4.149.55"'
$ws.Range('D15').Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Formula = '="  +2.01%  "'
$ws.Range('E15').Copy()
$ws.Range('E15').PasteSpecial(-4163)
$ws.Range('E16').Formula = '="  +0.58%  "'
$ws.Range('E16').Copy()
$ws.Range('E16').PasteSpecial(-4163)
$ws.Range('D17').Formula = '="612.34"'
$ws.Range('D17').Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Formula = '="  +0.39%  "'
$ws.Range('E17').Copy()
$ws.Range('E17').PasteSpecial(-4163)
$ws.Range('D18').Formula = '="3.575.15"'
$ws.Range('D18').Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('E18').Formula = '="  +2.06%  "'
$ws.Range('E18').Copy()
$ws.Range('E18').PasteSpecial(-4163)
$ws.Range('D19').Formula = '="70.738.54"'
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Formula = '="  +2.73%  "'
$ws.Range('E19').Copy()
$ws.Range('E19').PasteSpecial(-4163)
$ws.Range('E20').Formula = '="  -0.73%  "'
$ws.Range('E20').Copy()
$ws.Range('E20').PasteSpecial(-4163)
$ws.Range('D21').Formula = '="17.44"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Formula = '="  +0.90%  "'
$ws.Range('E21').Copy()
$ws.Range('E21').PasteSpecial(-4163)
$ws.Range('D22').Formula = '="0.885"'
$ws.Range('D22').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Formula = '="  +0.94%  "'
$ws.Range('E22').Copy()
$ws.Range('E22').PasteSpecial(-4163)
$ws.Range('E23').Formula = '="  -16.85%  "'
$ws.Range('E23').Copy()
$ws.Range('E23').PasteSpecial(-4163)
$ws.Range('D24').Formula = '="15.88"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Formula = '="  +1.74%  "'
$ws.Range('E24').Copy()
$ws.Range('E24').PasteSpecial(-4163)
$ws.Range('D25').Formula = '="97.13"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Formula = '="  +1.45%  "'
$ws.Range('E25').Copy()
$ws.Range('E25').PasteSpecial(-4163)
$ws.Range('E26').Formula = '="  -0.80%  "'
$ws.Range('E26').Copy()
$ws.Range('E26').PasteSpecial(-4163)
$ws.Range('E27').Formula = '="  -0.02%  "'
$ws.Range('E27').Copy()
$ws.Range('E27').PasteSpecial(-4163)
$ws.Range('E28').Formula = '="  +1.72%  "'
$ws.Range('E28').Copy()
$ws.Range('E28').PasteSpecial(-4163)
$ws.Range('E29').Formula = '="  +5.16%  "'
$ws.Range('E29').Copy()
$ws.Range('E29').PasteSpecial(-4163)
$ws.Range('D30').Formula = '="9.17"'
$ws.Range('D30').Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Formula = '="  -0.02%  "'
$ws.Range('E30').Copy()
$ws.Range('E30').PasteSpecial(-4163)
$ws.Range('D31').Formula = '="8.35"'
$ws.Range('D31').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Formula = '="  -0.62%  "'
$ws.Range('E31').Copy()
$ws.Range('E31').PasteSpecial(-4163)
$ws.Range('E32').Formula = '="  -1.15%  "'
$ws.Range('E32').Copy()
$ws.Range('E32').PasteSpecial(-4163)
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').Formula = '="667.43"'
$ws.Range('D33').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Formula = '="  +8.10%  "'
$ws.Range('E33').Copy()
$ws.Range('E33').PasteSpecial(-4163)
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Formula = '="7.18"'
$ws.Range('D34').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Formula = '="  +4.61%  "'
$ws.Range('E34').Copy()
$ws.Range('E34').PasteSpecial(-4163)
$ws.Range('E35').Formula = '="  -0.71%  "'
$ws.Range('E35').Copy()
$ws.Range('E35').PasteSpecial(-4163)
$ws.Range('D36').Formula = '="3.68"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Formula = '="  +7.81%  "'
$ws.Range('E36').Copy()
$ws.Range('E36').PasteSpecial(-4163)
$ws.Range('E38').Formula = '="  +1.63%  "'
$ws.Range('E38').Copy()
$ws.Range('E38').PasteSpecial(-4163)
$ws.Range('E39').Formula = '="  +8.70%  "'
$ws.Range('E39').Copy()
$ws.Range('E39').PasteSpecial(-4163)
$ws.Range('D40').Formula = '="57.34"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Formula = '="  +0.55%  "'
$ws.Range('E40').Copy()
$ws.Range('E40').PasteSpecial(-4163)
$ws.Range('E41').Formula = '="  -0.01%  "'
$ws.Range('E41').Copy()
$ws.Range('E41').PasteSpecial(-4163)
$ws.Range('E42').Formula = '="  +6.18%  "'
$ws.Range('E42').Copy()
$ws.Range('E42').PasteSpecial(-4163)
$ws.Range('D43').Formula = '="3.393.82"'
$ws.Range('D43').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Formula = '="  +0.77%  "'
$ws.Range('E43').Copy()
$ws.Range('E43').PasteSpecial(-4163)
$ws.Range('D44').Formula = '="0.322"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Formula = '="  -0.33%  "'
$ws.Range('E44').Copy()
$ws.Range('E44').PasteSpecial(-4163)
$ws.Range('D45').Formula = '="0.0₃0715"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Formula = '="  +4.09%  "'
$ws.Range('E45').Copy()
$ws.Range('E45').PasteSpecial(-4163)
$ws.Range('D46').Formula = '="32.92"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Formula = '="  +1.27%  "'
$ws.Range('E46').Copy()
$ws.Range('E46').PasteSpecial(-4163)
$ws.Range('D47').Formula = '="2.95"'
$ws.Range('D47').Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Formula = '="  +8.05%  "'
$ws.Range('E47').Copy()
$ws.Range('E47').PasteSpecial(-4163)
$ws.Range('D48').Formula = '="2.66"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Formula = '="  +6.26%  "'
$ws.Range('E48').Copy()
$ws.Range('E48').PasteSpecial(-4163)
$ws.Range('E49').Formula = '="  +1.45%  "'
$ws.Range('E49').Copy()
$ws.Range('E49').PasteSpecial(-4163)
$ws.Range('D50').Formula = '="132.30"'
$ws.Range('D50').Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Formula = '="  -0.52%  "'
$ws.Range('E50').Copy()
$ws.Range('E50').PasteSpecial(-4163)
$ws.Range('E51').Formula = '="  -0.05%  "'
$ws.Range('E51').Copy()
$ws.Range('E51').PasteSpecial(-4163)
